$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old extra data rows (3-12); only the header + one data row remain
$ws.Range("A3:F12").ClearContents()

# Re-arrange / relabel the header row (row 1) for the new column layout:
# A=total, B=dsMathSkills, C=AWSml, D=CV (new), E=matlab (new)
$ws.Range("A1").Value = "c_total_hours"
$ws.Range("B1").Value = "c_dsMathSkills_hours"
$ws.Range("C1").Value = "c_AWSml_hours"
$ws.Range("D1").Value = "c_CV_hours"
$ws.Range("E1").Value = "c_matlab_hours"

# Update the data row (row 2): total is now a SUM formula over the other 4 columns
$ws.Range("A2").Formula = "=SUM(B2:E2)"
$ws.Range("B2").Value = 13
$ws.Range("C2").Value = 8
$ws.Range("D2").Value = 10
$ws.Range("E2").Value = 35

# Resize columns B:E (bestFit-style widths) to fit the new header text
$ws.Columns.Item(2).ColumnWidth = 17.833333333333336
$ws.Columns.Item(3).ColumnWidth = 13.5
$ws.Columns.Item(4).ColumnWidth = 9.833333333333332
$ws.Columns.Item(5).ColumnWidth = 13.333333333333332

# Drop the now-unused 6th (F) column entirely
$ws.Columns.Item(6).Delete()

# Update the active selection to match the saved view
$ws.Range("G7").Select()
